$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B4").Value = 6.37
$ws.Range("D4").Value = -7.816
$ws.Range("E4").Value = 12.504
$ws.Range("D5").Value = -8.206999999999999
$ws.Range("B7").Value = 7.374
$ws.Range("D8").Value = -7.896000000000001
$ws.Range("E9").Value = 12.946
$ws.Range("B16").Value = 6.427
$ws.Range("D16").Value = -7.917
$ws.Range("E18").Value = 13.19
